$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day1")

$ws.Range("B2").Value = 90
$ws.Range("B5").Value = 10
